$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Corrected "web planform" values for rows 2-8, columns B:F
$data = @{
    2 = @(0.0871, 0.0871, 0.1059999987483025, 0.3151000142097473, -0.06889999657869339)
    3 = @(0.0726, 0.0726, 0.0513, 0.4611000120639801, -0.06889999657869339)
    4 = @(0.0358, 0.0358, 0.022, 0.3772000074386597, -0.06889999657869339)
    5 = @(-0.0016, -0.0016, -0.0104, 0.2377000004053116, -0.06889999657869339)
    6 = @(-0.06, -0.06, -0.0689, 0.2085999995470047, -0.06889999657869339)
    7 = @(-0.0613, -0.0613, -0.06889999657869339, 0.09390000253915787, -0.06889999657869339)
    8 = @(0.0764, 0.0764, 0.0703, 0.2233999967575073, -0.06889999657869339)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    $ws.Range("F$row").Value = $values[4]
}
